$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H1 header text (value only, style already set)
$ws.Range("H1").Value = "Contrast ratio"

# Add new I1 header; copy formatting (bold/centered) from H1, then set text
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Contrast ratio 95% CI"

# Row 2: H2 becomes numeric value, I2 gets new CI text
$ws.Range("H2").Value = 2.48
$ws.Range("I2").Value = "1.42-4.34"

# Row 3: H3 becomes numeric value, I3 gets new CI text
$ws.Range("H3").Value = 4.54
$ws.Range("I3").Value = "2.3-8.98"

# Row 4: H4 becomes numeric value, I4 gets new CI text
$ws.Range("H4").Value = 1.83
$ws.Range("I4").Value = "0.85-3.93"
